$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 1118 (i.e. directly
# after the existing row 1117), pushing the old rows 1118-1177 down to
# 1120-1179.
$ws.Range("A1118:A1119").EntireRow.Insert()

# Populate the two newly inserted rows (new row 1117 stays in place, the two
# new rows are 1118 and 1119) with the new daily price records for Naranja.

# New row 1118: Fukumoto / Primera
$ws.Cells.Item(1118, 1).Value2  = 9
$ws.Cells.Item(1118, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1118, 3).Value2  = "Metropolitana"
$ws.Cells.Item(1118, 4).Value2  = 45041
$ws.Cells.Item(1118, 5).Value2  = 13
$ws.Cells.Item(1118, 6).Value2  = "Fruta"
$ws.Cells.Item(1118, 7).Value2  = 100102
$ws.Cells.Item(1118, 8).Value2  = "Cítricos"
$ws.Cells.Item(1118, 9).Value2  = 100102005
$ws.Cells.Item(1118, 10).Value2 = "Naranja"
$ws.Cells.Item(1118, 11).Value2 = "Fukumoto"
$ws.Cells.Item(1118, 12).Value2 = "Primera"
$ws.Cells.Item(1118, 13).Value2 = 660
$ws.Cells.Item(1118, 14).Value2 = 10500
$ws.Cells.Item(1118, 15).Value2 = 12000
$ws.Cells.Item(1118, 16).Value2 = 11364
$ws.Cells.Item(1118, 17).Value2 = "$/caja 15 kilos granel"
$ws.Cells.Item(1118, 18).Value2 = "Provincia de Melipilla"
$ws.Cells.Item(1118, 19).Value2 = 758
$ws.Cells.Item(1118, 20).Value2 = 15

# New row 1119: Valencia / Primera
$ws.Cells.Item(1119, 1).Value2  = 9
$ws.Cells.Item(1119, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1119, 3).Value2  = "Metropolitana"
$ws.Cells.Item(1119, 4).Value2  = 45041
$ws.Cells.Item(1119, 5).Value2  = 13
$ws.Cells.Item(1119, 6).Value2  = "Fruta"
$ws.Cells.Item(1119, 7).Value2  = 100102
$ws.Cells.Item(1119, 8).Value2  = "Cítricos"
$ws.Cells.Item(1119, 9).Value2  = 100102005
$ws.Cells.Item(1119, 10).Value2 = "Naranja"
$ws.Cells.Item(1119, 11).Value2 = "Valencia"
$ws.Cells.Item(1119, 12).Value2 = "Primera"
$ws.Cells.Item(1119, 13).Value2 = 450
$ws.Cells.Item(1119, 14).Value2 = 11500
$ws.Cells.Item(1119, 15).Value2 = 12000
$ws.Cells.Item(1119, 16).Value2 = 11722
$ws.Cells.Item(1119, 17).Value2 = "$/caja 15 kilos granel"
$ws.Cells.Item(1119, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(1119, 19).Value2 = 781
$ws.Cells.Item(1119, 20).Value2 = 15

# Make sure column D keeps its date number format on the two new rows (the
# Insert already carries the style across, this just guarantees it).
$ws.Cells.Item(1118, 4).NumberFormat = $ws.Cells.Item(1117, 4).NumberFormat
$ws.Cells.Item(1119, 4).NumberFormat = $ws.Cells.Item(1117, 4).NumberFormat
